$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet
# contain the same event table and received the identical update.
foreach ($sheetName in "展览", "全部类型") {
    $ws = $wb.Worksheets.Item($sheetName)

    # The first data row (row 2, the "2024-03-16 OrangeOrange" event) was
    # removed from the source feed. Every subsequent event (rows 3..25)
    # moves up by one row. Column A holds a fixed running index (0,1,2,...)
    # per row position and must stay untouched, so only columns B:I are
    # shifted. Using Copy/Paste (rather than a plain Value assignment)
    # preserves the original cell types -- in particular it keeps the
    # date-like strings in column B as plain text instead of letting Excel
    # reinterpret them as date serials.
    $src = $ws.Range("B3:I25")
    $dst = $ws.Range("B2:I24")
    $src.Copy($dst)

    # The last row (25) is now a duplicate of row 24 and must be removed,
    # shrinking the sheet's used range back down to A1:I24.
    $ws.Rows.Item(25).Delete()
}

$excel.CutCopyMode = $false
